$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 203 ("Fruta/hortaliza, semanal").
# This pushes the existing rows 203-278 down to 204-279 (all of their cell
# contents, including formatting, travel with them), and we populate the
# freshly-inserted row 203 with the new record's data.
$ws.Rows(203).Insert()

$ws.Cells.Item(203, 1).Value = 10
$ws.Cells.Item(203, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(203, 3).Value = "La Araucanía"
$ws.Cells.Item(203, 4).Value = 44795
$ws.Cells.Item(203, 5).Value = 9
$ws.Cells.Item(203, 6).Value = 100112039
$ws.Cells.Item(203, 7).Value = "Ciboulette"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 40
$ws.Cells.Item(203, 11).Value = 8000
$ws.Cells.Item(203, 12).Value = 8000
$ws.Cells.Item(203, 13).Value = 8000
$ws.Cells.Item(203, 14).Value = "$/docena de atados"
$ws.Cells.Item(203, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(203, 16).Value = 2667
$ws.Cells.Item(203, 17).Value = 3
$ws.Cells.Item(203, 18).Value = "Hortaliza"
